$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$arr = New-Object 'object[,]' 89,2
$arr[0,0] = '253 Hoddle Street Apartment Complex Collingwood'; $arr[0,1] = 5
$arr[1,0] = '3535 Opal Meadow Heights Aged Care Community Meadow Heights'; $arr[1,1] = 27
$arr[2,0] = 'ABD Group 512 Melbourne Road Construction Site Spotswood'; $arr[2,1] = 5
$arr[3,0] = 'Acquire BPO Southbank'; $arr[3,1] = 7
$arr[4,0] = 'Al Haj Halal Meats Glenroy'; $arr[4,1] = 37
$arr[5,0] = 'Al-Taqwa College Truganina'; $arr[5,1] = 6
$arr[6,0] = 'Amiga Montessori Craigieburn'; $arr[6,1] = 25
$arr[7,0] = 'Best&Less Fountain Gate Narre Warren'; $arr[7,1] = 5
$arr[8,0] = 'Budget Car and Truck Rentals Campbellfield'; $arr[8,1] = 7
$arr[9,0] = 'CS Square Caroline Springs'; $arr[9,1] = 7
$arr[10,0] = 'Cannie Road Construction Site Cannie'; $arr[10,1] = 8
$arr[11,0] = 'Caroline Springs Police Station'; $arr[11,1] = 12
$arr[12,0] = 'Cedars Medical Clinic Coburg'; $arr[12,1] = 42
$arr[13,0] = 'Chemist Warehouse Fillo Drive Somerton'; $arr[13,1] = 5
$arr[14,0] = 'City of Hobsons Bay Community'; $arr[14,1] = 5
$arr[15,0] = 'City of Moreland Community'; $arr[15,1] = 7
$arr[16,0] = 'City of Wyndham Community'; $arr[16,1] = 7
$arr[17,0] = 'Classy Cabinets and Kitchens Craigieburn'; $arr[17,1] = 10
$arr[18,0] = 'Coles Aurora Village Epping'; $arr[18,1] = 6
$arr[19,0] = 'Coles Broadmeadows Central Shopping Centre'; $arr[19,1] = 9
$arr[20,0] = 'Coles Campbellfield Plaza Campbellfield'; $arr[20,1] = 8
$arr[21,0] = 'Coles Coburg North Village'; $arr[21,1] = 27
$arr[22,0] = 'Coles Greenvale Shopping Centre'; $arr[22,1] = 7
$arr[23,0] = 'Coles Pakenham Place Shopping Centre'; $arr[23,1] = 10
$arr[24,0] = 'Coles Roxburgh Village Roxburgh Park'; $arr[24,1] = 9
$arr[25,0] = 'Community Kids Meadow Heights'; $arr[25,1] = 14
$arr[26,0] = 'Construction Site Olea Apartment Caulfield North'; $arr[26,1] = 12
$arr[27,0] = 'Costco Wholesale Epping'; $arr[27,1] = 18
$arr[28,0] = 'Crusader Caravans Epping'; $arr[28,1] = 13
$arr[29,0] = 'DRC Laverton Automotive Repairs Laverton North'; $arr[29,1] = 5
$arr[30,0] = 'Direct Freight Express Campbellfield'; $arr[30,1] = 13
$arr[31,0] = 'Fitzroy Community School Fitzroy North'; $arr[31,1] = 41
$arr[32,0] = 'Fonterra Manufacturing Workplace Campbellfield'; $arr[32,1] = 7
$arr[33,0] = 'Glenroy West Primary School'; $arr[33,1] = 5
$arr[34,0] = 'Gumboots Early Learning Centre South Morang'; $arr[34,1] = 5
$arr[35,0] = 'Hamilton Marino 236 Jasper Road McKinnon'; $arr[35,1] = 12
$arr[36,0] = 'Health Care Providers Association South Melbourne'; $arr[36,1] = 13
$arr[37,0] = 'IGA Meadow Heights Shopping Centre Meadow Heights'; $arr[37,1] = 7
$arr[38,0] = 'ISS Factory Level 1 Terminal 2 Melbourne Airport Tullamarine'; $arr[38,1] = 7
$arr[39,0] = 'Ibis Kingsgate Hotel Melbourne'; $arr[39,1] = 6
$arr[40,0] = 'Ilim College Glenroy Campus Hadfield'; $arr[40,1] = 19
$arr[41,0] = 'Ilim Learning Sanctuary Glenroy'; $arr[41,1] = 12
$arr[42,0] = 'Industrial Galvanizers Valmont Coatings Campbellfield'; $arr[42,1] = 18
$arr[43,0] = 'Islamic College of Melbourne Tarneit'; $arr[43,1] = 9
$arr[44,0] = 'KFC Fawkner'; $arr[44,1] = 8
$arr[45,0] = 'Kasr Sweets Coolaroo'; $arr[45,1] = 6
$arr[46,0] = 'Kids House Early Learning Cheltenham'; $arr[46,1] = 12
$arr[47,0] = 'Learning Nest Early Learning Centre Meadow Heights'; $arr[47,1] = 6
$arr[48,0] = 'Level Crossing Removal Project Lilydale Construction Site John Street'; $arr[48,1] = 9
$arr[49,0] = 'Lineage Logistics Laverton North'; $arr[49,1] = 8
$arr[50,0] = 'Linfox Somerton National Distribution Centre Somerton'; $arr[50,1] = 9
$arr[51,0] = 'McDonalds Thomastown II'; $arr[51,1] = 7
$arr[52,0] = 'Melbourne Metropolitan Remand Centre Ravenhall'; $arr[52,1] = 11
$arr[53,0] = 'Melbourne Truck Repairs Campbellfield'; $arr[53,1] = 7
$arr[54,0] = 'Melbourne West Police Station Docklands'; $arr[54,1] = 7
$arr[55,0] = 'Melbourne Youth Justice Centre Parkville'; $arr[55,1] = 5
$arr[56,0] = 'Melton Police Station Melton'; $arr[56,1] = 5
$arr[57,0] = 'Mercy Hospital for Women Heidelberg'; $arr[57,1] = 5
$arr[58,0] = 'Mernda YMCA Early Learning Centre Mernda'; $arr[58,1] = 5
$arr[59,0] = 'Montessori Beginnings Greenvale'; $arr[59,1] = 5
$arr[60,0] = 'MyCentre Childcare Broadmeadows'; $arr[60,1] = 14
$arr[61,0] = 'National Gallery of Victoria Melbourne'; $arr[61,1] = 9
$arr[62,0] = 'Newbury Child and Community Centre Craigieburn'; $arr[62,1] = 7
$arr[63,0] = 'Nido Early School Moonee Ponds'; $arr[63,1] = 14
$arr[64,0] = 'Nido Early School Moonee Ponds'; $arr[64,1] = 14
$arr[65,0] = 'Nino Early Learning Adventures Lalor'; $arr[65,1] = 5
$arr[66,0] = 'North Geelong House Party'; $arr[66,1] = 7
$arr[67,0] = 'Northern Health Northern Hospital Epping Emergency Department Tier 1B'; $arr[67,1] = 44
$arr[68,0] = 'OnQ Plumbing and Excavations Craigieburn'; $arr[68,1] = 18
$arr[69,0] = 'Oporto Coolaroo'; $arr[69,1] = 13
$arr[70,0] = 'Paisley Park Early Learning Centre Bundoora'; $arr[70,1] = 9
$arr[71,0] = 'Panorama Construction Site Whitehorse Rd Box Hill'; $arr[71,1] = 20
$arr[72,0] = 'People First Healthcare Home Residence Disability Support Taylors Lakes'; $arr[72,1] = 5
$arr[73,0] = 'Ramsay Health Care Warringal Private Hospital Heidelberg'; $arr[73,1] = 9
$arr[74,0] = 'Richmond Quarter 261-271 Bridge Road Construction Site Richmond'; $arr[74,1] = 9
$arr[75,0] = 'Salta Drive Construction Site Rangedale Drainage Altona North'; $arr[75,1] = 7
$arr[76,0] = 'St Vincents Hospital Emergency Department Melbourne'; $arr[76,1] = 5
$arr[77,0] = 'Tek Foods Somerton'; $arr[77,1] = 13
$arr[78,0] = 'The Homestead Child and Family Centre Roxburgh Park'; $arr[78,1] = 13
$arr[79,0] = 'The Royal Children''s Hospital Melbourne Emergency Department Parkville Tier 1B'; $arr[79,1] = 10
$arr[80,0] = 'ThorwestenCabinets Pakenham'; $arr[80,1] = 13
$arr[81,0] = 'Total Window Concepts Hoppers Crossing'; $arr[81,1] = 6
$arr[82,0] = 'Unilodge College Square Student Accommodation 570 Lygon Street Carlton'; $arr[82,1] = 14
$arr[83,0] = 'Werribee Mercy Hospital Emergency Department'; $arr[83,1] = 8
$arr[84,0] = 'Western Health Footscray Hospital Emergency Department'; $arr[84,1] = 5
$arr[85,0] = 'Western Health Sunshine Hospital Emergency Department'; $arr[85,1] = 7
$arr[86,0] = 'Woodlands Long Day Care and Kindergarten Roxburgh Park'; $arr[86,1] = 5
$arr[87,0] = 'Woolworths Greenvale Lakes Roxburgh Park'; $arr[87,1] = 6
$arr[88,0] = 'Yara Childcare Centre Truganina'; $arr[88,1] = 10

$ws.Range("A2:B90").Value2 = $arr

# Remove now-extra trailing rows (91,92) left over from the old 92-row table
$ws.Range("A91:B92").ClearContents()